$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 6.820912598721751
$ws.Range("D2").Value = 5.70151577176208
$ws.Range("E2").Value = 11.08340569134277
$ws.Range("F2").Value = 48.71888702128586
$ws.Range("G2").Value = 63.53804066105725
$ws.Range("H2").Value = 22.46658599470801
$ws.Range("I2").Value = 37.36823565528921
$ws.Range("J2").Value = 9.865134272908122
$ws.Range("L2").Value = 9.465304517862631
$ws.Range("M2").Value = 30.84481856938868
$ws.Range("N2").Value = 17.52634659968407
# Row 3
$ws.Range("C3").Value = 6.826459800416182
$ws.Range("D3").Value = 5.666067956587002
$ws.Range("E3").Value = 11.10840268998274
$ws.Range("F3").Value = 48.58106236823549
$ws.Range("G3").Value = 63.07106968057143
$ws.Range("H3").Value = 22.45622469797813
$ws.Range("I3").Value = 37.29707625788979
$ws.Range("J3").Value = 9.899902076866796
$ws.Range("L3").Value = 9.488614189692694
$ws.Range("M3").Value = 30.21168592982617
$ws.Range("N3").Value = 17.3076605741715
# Row 4
$ws.Range("C4").Value = 6.830163008576514
$ws.Range("D4").Value = 5.643824383140236
$ws.Range("E4").Value = 11.12459823832196
$ws.Range("F4").Value = 48.51139695660602
$ws.Range("G4").Value = 62.80540223011759
$ws.Range("H4").Value = 22.45579872348409
$ws.Range("I4").Value = 37.2643849442028
$ws.Range("J4").Value = 9.92224022619502
$ws.Range("L4").Value = 9.503636415162136
$ws.Range("M4").Value = 29.81874467176163
$ws.Range("N4").Value = 17.17419911573255
# Row 5
$ws.Range("C5").Value = 6.831746882843774
$ws.Range("D5").Value = 5.634640923975969
$ws.Range("E5").Value = 11.13141171062487
$ws.Range("F5").Value = 48.48677523525924
$ws.Range("G5").Value = 62.70251844644347
$ws.Range("H5").Value = 22.45711421108927
$ws.Range("I5").Value = 37.25382830052191
$ws.Range("J5").Value = 9.931593192214077
$ws.Range("L5").Value = 9.509937218758703
$ws.Range("M5").Value = 29.65776032650312
$ws.Range("N5").Value = 17.12007578858776
# Row 6
$ws.Range("C6").Value = 6.832014402483477
$ws.Range("D6").Value = 5.633108866592345
$ws.Range("E6").Value = 11.13255600426192
$ws.Range("F6").Value = 48.48291444191224
$ws.Range("G6").Value = 62.68576169533245
$ws.Range("H6").Value = 22.45742245916151
$ws.Range("I6").Value = 37.25224231364644
$ws.Range("J6").Value = 9.933161370179585
$ws.Range("L6").Value = 9.510994298427688
$ws.Range("M6").Value = 29.63098309944927
$ws.Range("N6").Value = 17.11110627571253
# Row 7
$ws.Range("C7").Value = 6.83018406637613
$ws.Range("D7").Value = 5.643701012018192
$ws.Range("E7").Value = 11.12468926131151
$ws.Range("F7").Value = 48.51104964070982
$ws.Range("G7").Value = 62.80399282675431
$ws.Range("H7").Value = 22.45581044081066
$ws.Range("I7").Value = 37.26423138015738
$ws.Range("J7").Value = 9.922365350167874
$ws.Range("L7").Value = 9.503720663942525
$ws.Range("M7").Value = 29.81657678727797
$ws.Range("N7").Value = 17.17346804565996
# Row 8
$ws.Range("C8").Value = 6.822763620352496
$ws.Range("D8").Value = 5.689393443770365
$ws.Range("E8").Value = 11.09184916181768
$ws.Range("F8").Value = 48.66826513496307
$ws.Range("G8").Value = 63.37271275315851
$ws.Range("H8").Value = 22.46177958391948
$ws.Range("I8").Value = 37.34141696631151
$ws.Range("J8").Value = 9.876917162274077
$ws.Range("L8").Value = 9.473194691286066
$ws.Range("M8").Value = 30.62750139315782
$ws.Range("N8").Value = 17.45081035075837
# Row 9
$ws.Range("C9").Value = 6.810568131364006
$ws.Range("D9").Value = 5.775145070943776
$ws.Range("E9").Value = 11.03414536228465
$ws.Range("F9").Value = 49.09479461685689
$ws.Range("G9").Value = 64.65082242231557
$ws.Range("H9").Value = 22.52069972844828
$ws.Range("I9").Value = 37.57998035778463
$ws.Range("J9").Value = 9.795611500614854
$ws.Range("L9").Value = 9.418939559694675
$ws.Range("M9").Value = 32.17614066687347
$ws.Range("N9").Value = 17.99854515740062
# Row 10
$ws.Range("C10").Value = 6.803041603652122
$ws.Range("D10").Value = 5.83572929195931
$ws.Range("E10").Value = 10.99579376500471
$ws.Range("F10").Value = 49.47931608217024
$ws.Range("G10").Value = 65.68278685169965
$ws.Range("H10").Value = 22.59286754131061
$ws.Range("I10").Value = 37.80812097845305
$ws.Range("J10").Value = 9.740583108451869
$ws.Range("L10").Value = 9.382457431863452
$ws.Range("M10").Value = 33.27807037803947
$ws.Range("N10").Value = 18.39990522103325
# Row 11
$ws.Range("C11").Value = 6.799928329260833
$ws.Range("D11").Value = 5.862751079188586
$ws.Range("E11").Value = 10.97921641675754
$ws.Range("F11").Value = 49.66940076567146
$ws.Range("G11").Value = 66.17089854113502
$ws.Range("H11").Value = 22.63196735109986
$ws.Range("I11").Value = 37.9232647356522
$ws.Range("J11").Value = 9.716558575998475
$ws.Range("L11").Value = 9.366586234327405
$ws.Range("M11").Value = 33.76954615582107
$ws.Range("N11").Value = 18.58155763467967
# Row 12
$ws.Range("C12").Value = 6.798794037504333
$ws.Range("D12").Value = 5.872904810331963
$ws.Range("E12").Value = 10.97306334397866
$ws.Range("F12").Value = 49.74352947077031
$ws.Range("G12").Value = 66.35827794765703
$ws.Range("H12").Value = 22.64767292315524
$ws.Range("I12").Value = 37.96848524355649
$ws.Range("J12").Value = 9.707605116611363
$ws.Range("L12").Value = 9.360679808812925
$ws.Range("M12").Value = 33.95409034013471
$ws.Range("N12").Value = 18.65015397479865
# Row 13
$ws.Range("C13").Value = 6.799036342607844
$ws.Range("D13").Value = 5.870721551408784
$ws.Range("E13").Value = 10.9743829938831
$ws.Range("F13").Value = 49.72746957253627
$ws.Range("G13").Value = 66.31781170077053
$ws.Range("H13").Value = 22.64425049876795
$ws.Range("I13").Value = 37.95867451381095
$ws.Range("J13").Value = 9.70952700868461
$ws.Range("L13").Value = 9.36194726213578
$ws.Range("M13").Value = 33.91441754689875
$ws.Range("N13").Value = 18.63538994658063
# Row 14
$ws.Range("C14").Value = 6.79983411591436
$ws.Range("D14").Value = 5.8635880208317
$ws.Range("E14").Value = 10.97870770955111
$ws.Range("F14").Value = 49.67545653742591
$ws.Range("G14").Value = 66.18626424399258
$ws.Range("H14").Value = 22.63324145162224
$ws.Range("I14").Value = 37.92695272549113
$ws.Range("J14").Value = 9.715819086271329
$ws.Range("L14").Value = 9.366098234993235
$ws.Range("M14").Value = 33.78476099049762
$ws.Range("N14").Value = 18.58720526045361
# Row 15
$ws.Range("C15").Value = 6.800328587593511
$ws.Range("D15").Value = 5.859208218778912
$ws.Range("E15").Value = 10.98137290825112
$ws.Range("F15").Value = 49.64387570771344
$ws.Range("G15").Value = 66.10601437155439
$ws.Range("H15").Value = 22.6266151219959
$ws.Range("I15").Value = 37.90773242032546
$ws.Range("J15").Value = 9.719691907225386
$ws.Range("L15").Value = 9.368654307771555
$ws.Range("M15").Value = 33.70513402084433
$ws.Range("N15").Value = 18.55766412045284
# Row 16
$ws.Range("C16").Value = 6.803251311075771
$ws.Range("D16").Value = 5.83395240058711
$ws.Range("E16").Value = 10.99689457284532
$ws.Range("F16").Value = 49.46719603240226
$ws.Range("G16").Value = 65.6512512048276
$ws.Range("H16").Value = 22.59043827715989
$ws.Range("I16").Value = 37.80082340718272
$ws.Range("J16").Value = 9.742173381135544
$ws.Range("L16").Value = 9.383509187034669
$ws.Range("M16").Value = 33.24574036652477
$ws.Range("N16").Value = 18.3880101314955
# Row 17
$ws.Range("C17").Value = 6.805123837008312
$ws.Range("D17").Value = 5.818319887109121
$ws.Range("E17").Value = 11.0066387872661
$ws.Range("F17").Value = 49.36267065349339
$ws.Range("G17").Value = 65.3769503907188
$ws.Range("H17").Value = 22.56984962479294
$ws.Range("I17").Value = 37.73813786121409
$ws.Range("J17").Value = 9.756222636875696
$ws.Range("L17").Value = 9.392807390347452
$ws.Range("M17").Value = 32.9612919574645
$ws.Range("N17").Value = 18.28365333846401
# Row 18
$ws.Range("C18").Value = 6.806230098993963
$ws.Range("D18").Value = 5.809277898266369
$ws.Range("E18").Value = 11.01232522744682
$ws.Range("F18").Value = 49.3039799644363
$ws.Range("G18").Value = 65.22094566092967
$ws.Range("H18").Value = 22.55859798008931
$ws.Range("I18").Value = 37.70315330859872
$ws.Range("J18").Value = 9.764398341403972
$ws.Range("L18").Value = 9.398223717022919
$ws.Range("M18").Value = 32.79677134956395
$ws.Range("N18").Value = 18.22354606344534
# Row 19
$ws.Range("C19").Value = 6.806609682120104
$ws.Range("D19").Value = 5.806207792729136
$ws.Range("E19").Value = 11.01426462882629
$ws.Range("F19").Value = 49.28435475995941
$ws.Range("G19").Value = 65.16843265407546
$ws.Range("H19").Value = 22.5548898215798
$ws.Range("I19").Value = 37.69149243395042
$ws.Range("J19").Value = 9.767182827517757
$ws.Range("L19").Value = 9.400069330546188
$ws.Range("M19").Value = 32.74091561271748
$ws.Range("N19").Value = 18.20318213984744
# Row 20
$ws.Range("C20").Value = 6.804921477934998
$ws.Range("D20").Value = 5.819989240143833
$ws.Range("E20").Value = 11.00559303434489
$ws.Range("F20").Value = 49.37364983017483
$ws.Range("G20").Value = 65.40596843274298
$ws.Range("H20").Value = 22.57198023084534
$ws.Range("I20").Value = 37.74470015508221
$ws.Range("J20").Value = 9.754717249507539
$ws.Range("L20").Value = 9.391810521581995
$ws.Range("M20").Value = 32.99166768962334
$ws.Range("N20").Value = 18.29477141738903
# Row 21
$ws.Range("C21").Value = 6.799598579277359
$ws.Range("D21").Value = 5.865685462518474
$ws.Range("E21").Value = 10.97743406366391
$ws.Range("F21").Value = 49.69067601934441
$ws.Range("G21").Value = 66.22483508905717
$ws.Range("H21").Value = 22.63645069249591
$ws.Range("I21").Value = 37.93622640880138
$ws.Range("J21").Value = 9.713967046332691
$ws.Range("L21").Value = 9.364876184992037
$ws.Range("M21").Value = 33.82288800025386
$ws.Range("N21").Value = 18.60136392110622
# Row 22
$ws.Range("C22").Value = 6.79637993083186
$ws.Range("D22").Value = 5.895090432554804
$ws.Range("E22").Value = 10.95975544747198
$ws.Range("F22").Value = 49.91037160819531
$ws.Range("G22").Value = 66.7747594578408
$ws.Range("H22").Value = 22.68382592866351
$ws.Range("I22").Value = 38.07082197592731
$ws.Range("J22").Value = 9.688174014819399
$ws.Range("L22").Value = 9.347876946864348
$ws.Range("M22").Value = 34.3569388567926
$ws.Range("N22").Value = 18.80059810624113
# Row 23
$ws.Range("C23").Value = 6.798073985285363
$ws.Range("D23").Value = 5.879439019178558
$ws.Range("E23").Value = 10.96912471028126
$ws.Range("F23").Value = 49.79198404281194
$ws.Range("G23").Value = 66.47995339843914
$ws.Range("H23").Value = 22.65806247085082
$ws.Range("I23").Value = 37.99812957175075
$ws.Range("J23").Value = 9.701863703959859
$ws.Range("L23").Value = 9.35689468693398
$ws.Range("M23").Value = 34.0727971287791
$ws.Range("N23").Value = 18.69438623465513
# Row 24
$ws.Range("C24").Value = 6.805012871893378
$ws.Range("D24").Value = 5.819234695746534
$ws.Range("E24").Value = 11.00606555643146
$ws.Range("F24").Value = 49.36868177397253
$ws.Range("G24").Value = 65.39284408736795
$ws.Range("H24").Value = 22.57101516139763
$ws.Range("I24").Value = 37.7417300572607
$ws.Range("J24").Value = 9.755397527953503
$ws.Range("L24").Value = 9.392260985800375
$ws.Range("M24").Value = 32.97793788462007
$ws.Range("N24").Value = 18.28974527574188
# Row 25
$ws.Range("C25").Value = 6.813615404275901
$ws.Range("D25").Value = 5.752366471136473
$ws.Range("E25").Value = 11.04904294505682
$ws.Range("F25").Value = 48.96682674856959
$ws.Range("G25").Value = 64.28824301662884
$ws.Range("H25").Value = 22.49969136253206
$ws.Range("I25").Value = 37.50612555573295
$ws.Range("J25").Value = 9.816775985088277
$ws.Range("L25").Value = 9.433020821161159
$ws.Range("M25").Value = 31.76276497906785
$ws.Range("N25").Value = 17.85031052613873
